# Applies corrected TSLA price/shares/ticker data to rows 2-44 of Sheet1.
# Each data row is: Row,Open,Close,High,Low,SharesOutstanding,SetTickerToTSLA(0/1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @(2, 12.57999992370606, 15.06999969482422, 15.9166669845581, 12.40333271026611, 3325819167, 1),
    @(3, 18.07399940490723, 17.74333381652832, 19.11000061035156, 16.71933364868164, 3325819167, 1),
    @(4, 16.50066757202148, 13.79533290863037, 16.6560001373291, 13.46666717529297, 3325819167, 1),
    @(5, 15.38133335113525, 12.74666690826416, 15.42533302307129, 12.16066741943359, 3325819167, 1),
    @(6, 16.32200050354004, 16.05066680908203, 17.95599937438965, 15.55000019073486, 3325819167, 1),
    @(7, 13.74266719818115, 15.65266704559326, 15.68533325195312, 13.73333263397217, 3325819167, 1),
    @(8, 14.15333271026611, 13.18200016021728, 14.37800025939941, 12.80000019073486, 3325819167, 1),
    @(9, 14.32400035858154, 16.79533386230469, 17.23066711425781, 14.06400012969971, 3325819167, 1),
    @(10, 19.12666702270508, 20.9379997253418, 20.98666763305664, 18.97200012207031, 3325819167, 1),
    @(11, 24.68266677856445, 21.56466674804688, 24.75666618347168, 20.20866775512696, 3325819167, 1),
    @(12, 22.83466720581055, 22.10199928283692, 24.20000076293945, 21.11066627502441, 3325819167, 1),
    @(13, 20.79999923706055, 23.62066650390625, 24.03333282470703, 20.3786678314209, 3325819167, 1),
    @(14, 17.08399963378906, 19.59333229064941, 20.63333320617676, 16.30599975585938, 3325819167, 1),
    @(15, 24.00466728210449, 19.8759994506836, 24.31866645812988, 19.07533264160156, 3325819167, 1),
    @(16, 20.38466644287109, 22.48800086975098, 23.14399909973145, 16.51799964904785, 3325819167, 1),
    @(17, 20.40666770935058, 20.4680004119873, 23.46666717529297, 18.61866760253906, 3325819167, 1),
    @(18, 18.84133338928223, 15.9126672744751, 19.74466705322266, 15.40866661071777, 3325819167, 1),
    @(19, 15.34733295440674, 16.10733222961426, 17.73800086975098, 14.81466674804688, 3325819167, 1),
    @(20, 16.10000038146973, 20.99466705322266, 22.7226676940918, 14.95199966430664, 3325819167, 1),
    @(21, 28.29999923706055, 43.37133407592773, 43.53333282470703, 28.11400032043457, 3325819167, 1),
    @(22, 33.59999847412109, 52.12533187866211, 57.98799896240234, 29.76000022888184, 3325819167, 1),
    @(23, 72.19999694824219, 95.38400268554688, 119.6660003662109, 72.03333282470703, 3325819167, 1),
    @(24, 146.9199981689453, 129.3466644287109, 155.3000030517578, 126.370002746582, 3325819167, 1),
    @(25, 239.8200073242188, 264.510009765625, 300.1333312988281, 239.0633392333984, 3325819167, 1),
    @(26, 229.4566650390625, 236.479995727539, 260.2633361816406, 219.8066711425781, 3325819167, 1),
    @(27, 227.9733276367188, 229.0666656494141, 233.3333282470703, 206.8200073242188, 3325819167, 1),
    @(28, 259.4666748046875, 371.3333435058594, 371.7366638183594, 254.5299987792969, 3325819167, 1),
    @(29, 382.5833435058594, 312.239990234375, 402.6666564941406, 264.0033264160156, 3325819167, 1),
    @(30, 360.3833312988281, 290.2533264160156, 384.2900085449219, 273.8999938964844, 3325819167, 1),
    @(31, 227, 297.1499938964844, 298.3200073242188, 216.1666717529297, 3325819167, 1),
    @(32, $null, $null, $null, $null, 3325819167, 0),
    @(33, 118.4700012207031, 173.2200012207031, 180.6799926757812, 101.8099975585938, 3325819167, 1),
    @(34, 199.9100036621093, 164.3099975585938, 202.6900024414062, 152.3699951171875, 3325819167, 1),
    @(35, 276.489990234375, 267.4299926757812, 299.2900085449219, 254.1199951171875, 3325819167, 1),
    @(36, 244.8099975585937, 200.8399963378907, 268.9400024414062, 194.0700073242188, 3325819167, 1),
    @(37, 250.0800018310547, 187.2899932861328, 251.25, 180.0599975585938, 3325819167, 1),
    @(38, 176.1699981689453, 183.2799987792969, 198.8699951171875, 138.8000030517578, 3325819167, 1),
    @(39, 201.0200042724609, 232.0700073242188, 271, 200.8500061035156, 3325819167, 1),
    @(40, 262.6700134277344, 249.8500061035156, 273.5400085449219, 212.1100006103516, 3325819167, 1),
    @(41, 390.1000061035156, 404.6000061035156, 439.739990234375, 373.0400085449219, 3325819167, 1),
    @(42, 263.7999877929688, 282.1600036621094, 294.8599853515625, 214.25, 3325819167, 1),
    @(43, 298.4599914550781, 308.2699890136719, 338, 288.7699890136719, 3325819167, 1),
    @(44, 443.7999877929688, 433.7200012207031, 470.75, 411.4500122070313, 3325819167, 1)
)

foreach ($r in $rowData) {
    $row = $r[0]
    if ($null -ne $r[1]) { $ws.Range("D$row").Value = [double]$r[1] }
    if ($null -ne $r[2]) { $ws.Range("E$row").Value = [double]$r[2] }
    if ($null -ne $r[3]) { $ws.Range("F$row").Value = [double]$r[3] }
    if ($null -ne $r[4]) { $ws.Range("G$row").Value = [double]$r[4] }
    if ($null -ne $r[5]) { $ws.Range("H$row").Value = [double]$r[5] }
    if ($r[6] -eq 1) { $ws.Range("I$row").Value = "TSLA" }
}
